# Auto-generated edit script: refreshes the market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all leve-profit
# sheets to match the latest scraped Market Board data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4437.375
$ws.Range("J29").Value = 5666.6665
$ws.Range("L29").Value = 16999.9995
$ws.Range("N29").Value = -17561.9995
$ws.Range("H31").Value = 819.75
$ws.Range("I31").Value = 407.14285
$ws.Range("K31").Value = 1221.42855
$ws.Range("M31").Value = -991.4285500000001
$ws.Range("H82").Value = 29385
$ws.Range("I82").Value = 15846.667
$ws.Range("K82").Value = 47540.001
$ws.Range("M82").Value = -47134.001
$ws.Range("H85").Value = 29385
$ws.Range("I85").Value = 15846.667
$ws.Range("K85").Value = 47540.001
$ws.Range("M85").Value = -46136.001
$ws.Range("H100").Value = 1370
$ws.Range("I100").Value = 1160.3334
$ws.Range("J100").Value = 1999
$ws.Range("K100").Value = 1160.3334
$ws.Range("L100").Value = 1999
$ws.Range("M100").Value = -619.3334
$ws.Range("N100").Value = -3081
$ws.Range("H132").Value = 3015.1667
$ws.Range("I132").Value = 3015.1667
$ws.Range("K132").Value = 9045.500100000001
$ws.Range("M132").Value = -6515.500100000001
$ws.Range("H138").Value = 2106.8086
$ws.Range("I138").Value = 1563.3077
$ws.Range("J138").Value = 2194.037
$ws.Range("K138").Value = 4689.9231
$ws.Range("L138").Value = 6582.110999999999
$ws.Range("M138").Value = 450.0769
$ws.Range("N138").Value = -16862.111

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 506.48276
$ws.Range("I2").Value = 402.34784
$ws.Range("J2").Value = 905.6667
$ws.Range("K2").Value = 402.34784
$ws.Range("L2").Value = 905.6667
$ws.Range("M2").Value = -289.34784
$ws.Range("N2").Value = -1131.6667
$ws.Range("H32").Value = 10213785
$ws.Range("I32").Value = 13896017
$ws.Range("J32").Value = 16834.154
$ws.Range("K32").Value = 13896017
$ws.Range("L32").Value = 16834.154
$ws.Range("M32").Value = -13895730
$ws.Range("N32").Value = -17408.154
$ws.Range("H74").Value = 9633695
$ws.Range("I74").Value = 11906628
$ws.Range("K74").Value = 11906628
$ws.Range("M74").Value = -11905754
$ws.Range("H77").Value = 9633695
$ws.Range("I77").Value = 11906628
$ws.Range("K77").Value = 59533140
$ws.Range("M77").Value = -59528772
$ws.Range("H94").Value = 47744.75
$ws.Range("J94").Value = 47744.75
$ws.Range("L94").Value = 47744.75
$ws.Range("N94").Value = -49546.75
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H116").Value = 506.48276
$ws.Range("I116").Value = 402.34784
$ws.Range("J116").Value = 905.6667
$ws.Range("K116").Value = 402.34784
$ws.Range("L116").Value = 905.6667
$ws.Range("M116").Value = 1891.65216
$ws.Range("N116").Value = -5493.6667
$ws.Range("H141").Value = 90999.2
$ws.Range("I141").Value = 90999.2
$ws.Range("K141").Value = 90999.2
$ws.Range("M141").Value = -85819.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 506.48276
$ws.Range("I3").Value = 402.34784
$ws.Range("J3").Value = 905.6667
$ws.Range("K3").Value = 402.34784
$ws.Range("L3").Value = 905.6667
$ws.Range("M3").Value = -288.34784
$ws.Range("N3").Value = -1133.6667
$ws.Range("H21").Value = 110000
$ws.Range("J21").Value = 110000
$ws.Range("L21").Value = 110000
$ws.Range("N21").Value = -110472
$ws.Range("H36").Value = 2166.75
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H80").Value = 1250.8889
$ws.Range("I80").Value = 493
$ws.Range("J80").Value = 1345.625
$ws.Range("K80").Value = 493
$ws.Range("L80").Value = 1345.625
$ws.Range("M80").Value = 505
$ws.Range("N80").Value = -3341.625
$ws.Range("H83").Value = 1250.8889
$ws.Range("I83").Value = 493
$ws.Range("J83").Value = 1345.625
$ws.Range("K83").Value = 2465
$ws.Range("L83").Value = 6728.125
$ws.Range("M83").Value = 2527
$ws.Range("N83").Value = -16712.125
$ws.Range("H134").Value = 626605.5
$ws.Range("I134").Value = 1711.6
$ws.Range("K134").Value = 5134.799999999999
$ws.Range("M134").Value = -2599.799999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 751.4286
$ws.Range("I16").Value = 398.3
$ws.Range("K16").Value = 398.3
$ws.Range("M16").Value = -111.3
$ws.Range("H31").Value = 422663.28
$ws.Range("I31").Value = 6142.2197
$ws.Range("J31").Value = 1561154.2
$ws.Range("K31").Value = 6142.2197
$ws.Range("L31").Value = 1561154.2
$ws.Range("M31").Value = -5847.2197
$ws.Range("N31").Value = -1561744.2
$ws.Range("H34").Value = 422663.28
$ws.Range("I34").Value = 6142.2197
$ws.Range("J34").Value = 1561154.2
$ws.Range("K34").Value = 6142.2197
$ws.Range("L34").Value = 1561154.2
$ws.Range("M34").Value = -5940.2197
$ws.Range("N34").Value = -1561558.2
$ws.Range("H50").Value = 59285.43
$ws.Range("J50").Value = 59285.43
$ws.Range("L50").Value = 59285.43
$ws.Range("N50").Value = -60535.43
$ws.Range("H105").Value = 1804.3334
$ws.Range("I105").Value = 1810.5
$ws.Range("J105").Value = 1801.25
$ws.Range("K105").Value = 1810.5
$ws.Range("L105").Value = 1801.25
$ws.Range("M105").Value = -63.5
$ws.Range("N105").Value = -5295.25
$ws.Range("H111").Value = 66747.25
$ws.Range("J111").Value = 66747.25
$ws.Range("L111").Value = 66747.25
$ws.Range("N111").Value = -74927.25
$ws.Range("H113").Value = 751.4286
$ws.Range("I113").Value = 398.3
$ws.Range("K113").Value = 398.3
$ws.Range("M113").Value = 1771.7
$ws.Range("H132").Value = 3137.7693
$ws.Range("I132").Value = 2708.5454
$ws.Range("J132").Value = 5498.5
$ws.Range("K132").Value = 8125.6362
$ws.Range("L132").Value = 16495.5
$ws.Range("M132").Value = -5595.6362
$ws.Range("N132").Value = -21555.5
$ws.Range("H134").Value = 3525.5557
$ws.Range("I134").Value = 2704.077
$ws.Range("J134").Value = 5661.4
$ws.Range("K134").Value = 8112.231000000001
$ws.Range("L134").Value = 16984.2
$ws.Range("M134").Value = -5577.231000000001
$ws.Range("N134").Value = -22054.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1199
$ws.Range("I12").Value = 2028
$ws.Range("J12").Value = 497.53845
$ws.Range("K12").Value = 6084
$ws.Range("L12").Value = 1492.61535
$ws.Range("M12").Value = -5911
$ws.Range("N12").Value = -1838.61535
$ws.Range("H37").Value = 84461.75
$ws.Range("J37").Value = 84461.75
$ws.Range("L37").Value = 253385.25
$ws.Range("N37").Value = -253609.25
$ws.Range("H46").Value = 1683.1428
$ws.Range("J46").Value = 2196.6
$ws.Range("L46").Value = 6589.799999999999
$ws.Range("N46").Value = -6771.799999999999
$ws.Range("H56").Value = 5956.8667
$ws.Range("I56").Value = 5956.8667
$ws.Range("K56").Value = 5956.8667
$ws.Range("M56").Value = -5426.8667
$ws.Range("H132").Value = 1749.375
$ws.Range("I132").Value = 999.5
$ws.Range("J132").Value = 2499.25
$ws.Range("K132").Value = 8995.5
$ws.Range("L132").Value = 22493.25
$ws.Range("M132").Value = -6465.5
$ws.Range("N132").Value = -27553.25
$ws.Range("H140").Value = 2257.92
$ws.Range("I140").Value = 2227
$ws.Range("K140").Value = 6681
$ws.Range("M140").Value = -1501

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 112994.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 112994.5
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 112994.5
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -115158.5
$ws.Range("H102").Value = 2109.0286
$ws.Range("I102").Value = 1631.1154
$ws.Range("J102").Value = 3489.6667
$ws.Range("K102").Value = 1631.1154
$ws.Range("L102").Value = 3489.6667
$ws.Range("M102").Value = -9.115399999999909
$ws.Range("N102").Value = -6733.6667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 62504996
$ws.Range("J44").Value = 62504996
$ws.Range("L44").Value = 62504996
$ws.Range("N44").Value = -62505908
$ws.Range("H58").Value = 7428.5713
$ws.Range("I58").Value = 5000
$ws.Range("J58").Value = 8400
$ws.Range("K58").Value = 5000
$ws.Range("L58").Value = 8400
$ws.Range("M58").Value = -4740
$ws.Range("N58").Value = -8920
$ws.Range("H96").Value = 96400
$ws.Range("J96").Value = 96400
$ws.Range("L96").Value = 96400
$ws.Range("N96").Value = -101892
$ws.Range("H109").Value = 98758.336
$ws.Range("J109").Value = 98758.336
$ws.Range("L109").Value = 98758.336
$ws.Range("N109").Value = -101532.336
$ws.Range("H123").Value = 53323.332
$ws.Range("J123").Value = 53323.332
$ws.Range("L123").Value = 53323.332
$ws.Range("N123").Value = -63123.332
$ws.Range("H136").Value = 93269.92999999999
$ws.Range("I136").Value = 12781.556
$ws.Range("J136").Value = 214002.5
$ws.Range("K136").Value = 38344.66800000001
$ws.Range("L136").Value = 642007.5
$ws.Range("M136").Value = -35794.66800000001
$ws.Range("N136").Value = -647107.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 61089.668
$ws.Range("I31").Value = 47500
$ws.Range("K31").Value = 47500
$ws.Range("M31").Value = -47152
$ws.Range("H107").Value = 642.5854
$ws.Range("I107").Value = 679.2
$ws.Range("K107").Value = 2037.6
$ws.Range("M107").Value = -117.6000000000001
$ws.Range("H109").Value = 104960
$ws.Range("J109").Value = 104960
$ws.Range("L109").Value = 104960
$ws.Range("N109").Value = -107734
$ws.Range("H114").Value = 64000.5
$ws.Range("J114").Value = 64000.5
$ws.Range("L114").Value = 64000.5
$ws.Range("N114").Value = -72678.5
$ws.Range("H132").Value = 530958.1
$ws.Range("I132").Value = 4262.5
$ws.Range("J132").Value = 3340001.2
$ws.Range("K132").Value = 12787.5
$ws.Range("L132").Value = 10020003.6
$ws.Range("M132").Value = -10257.5
$ws.Range("N132").Value = -10025063.6
